# Fruta / hortaliza, semanal
# Weekly data refresh for "Feria Lagunitas de Puerto Montt - Chirimoya":
# two new observations (dated 2021-11-23 / serial 44523) are inserted as the
# new rows 12-13, pushing the former rows 12-13 down to rows 14-15.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows before the current row 12 (shifts old 12,13 -> 14,15)
$ws.Rows("12:13").Insert()

# New row 12: Primera quality, 2021-11-23
$ws.Range("A12").Value = 4
$ws.Range("B12").Value = 'Feria Lagunitas de Puerto Montt'
$ws.Range("C12").Value = 'Los Lagos'
$ws.Range("D12").Value = 44523
$ws.Range("E12").Value = 10
$ws.Range("F12").Value = 'Fruta'
$ws.Range("G12").Value = 100107
$ws.Range("H12").Value = 'Otros'
$ws.Range("I12").Value = 100107002
$ws.Range("J12").Value = 'Chirimoya'
$ws.Range("K12").Value = 'Cultivar IV Región'
$ws.Range("L12").Value = 'Primera'
$ws.Range("M12").Value = 400
$ws.Range("N12").Value = 21000
$ws.Range("O12").Value = 22000
$ws.Range("P12").Value = 21500
$ws.Range("Q12").Value = '$/bandeja 8 kilos'
$ws.Range("R12").Value = 'Provincia de Limarí'
$ws.Range("S12").Value = 2688
$ws.Range("T12").Value = 8

# New row 13: Segunda quality, 2021-11-23
$ws.Range("A13").Value = 4
$ws.Range("B13").Value = 'Feria Lagunitas de Puerto Montt'
$ws.Range("C13").Value = 'Los Lagos'
$ws.Range("D13").Value = 44523
$ws.Range("E13").Value = 10
$ws.Range("F13").Value = 'Fruta'
$ws.Range("G13").Value = 100107
$ws.Range("H13").Value = 'Otros'
$ws.Range("I13").Value = 100107002
$ws.Range("J13").Value = 'Chirimoya'
$ws.Range("K13").Value = 'Cultivar IV Región'
$ws.Range("L13").Value = 'Segunda'
$ws.Range("M13").Value = 100
$ws.Range("N13").Value = 18000
$ws.Range("O13").Value = 18000
$ws.Range("P13").Value = 18000
$ws.Range("Q13").Value = '$/bandeja 8 kilos'
$ws.Range("R13").Value = 'Provincia de Limarí'
$ws.Range("S13").Value = 2250
$ws.Range("T13").Value = 8
